$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, shifting existing rows 5..95 down to 6..96
$ws.Rows("5:5").Insert()

# Populate the newly inserted row 5 with the new market-price record
$ws.Range("A5").Value = 10
$ws.Range("B5").Value = "Vega Modelo de Temuco"
$ws.Range("C5").Value = "La Araucanía"
$ws.Range("D5").Value = "2022-02-24"
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = 100112022
$ws.Range("G5").Value = "Arveja Verde"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 23000
$ws.Range("L5").Value = 23000
$ws.Range("M5").Value = 23000
$ws.Range("N5").Value = "$/saco 25 kilos"
$ws.Range("O5").Value = "Región de La Araucanía"
$ws.Range("P5").Value = 920
$ws.Range("Q5").Value = 25
$ws.Range("R5").Value = "Hortaliza"
